$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial value that was updated for all
# data rows (rows 2 through 117) from 45172 (2023-09-03) to 45175 (2023-09-06).
$ws.Range("C2:C117").Value = 45175
